$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added. It belongs chronologically
# before the existing row 74, so insert a new row at position 74 which
# pushes the former rows 74-86 down to 75-87 (preserving all of their data
# and formatting), then fill the freshly inserted row 74 with the new data.
$ws.Rows.Item(74).Insert()

$ws.Cells.Item(74, 1).Value = 11
$ws.Cells.Item(74, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(74, 3).Value = "Bíobío"
$ws.Cells.Item(74, 4).Value = 44943
$ws.Cells.Item(74, 5).Value = 8
$ws.Cells.Item(74, 6).Value = 100112037
$ws.Cells.Item(74, 7).Value = "Cebollín"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 450
$ws.Cells.Item(74, 11).Value = 2700
$ws.Cells.Item(74, 12).Value = 2800
$ws.Cells.Item(74, 13).Value = 2756
$ws.Cells.Item(74, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(74, 15).Value = "Región Metropolitana"
$ws.Cells.Item(74, 16).Value = 77
$ws.Cells.Item(74, 17).Value = 36
$ws.Cells.Item(74, 18).Value = "Hortaliza"
